$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.362.65'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.047.29'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('D5').Value = '228.85'
$ws.Range('E5').Value = '  -1.71%  '
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '56.95'
$ws.Range('E8').Value = '  -2.62%  '
$ws.Range('D9').Value = '0.386'
$ws.Range('E9').Value = '  -1.58%  '
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').Value = '  +0.53%  '
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').Value = '14.78'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '2.333.72'
$ws.Range('E13').Value = '  -2.03%  '
$ws.Range('D14').Value = '20.73'
$ws.Range('E14').Value = '  -1.36%  '
$ws.Range('D15').Value = '0.757'
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').Value = '5.31'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '2.048.72'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = '37.268.01'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '6.12'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').Value = '69.51'
$ws.Range('E20').Value = '  -2.70%  '
$ws.Range('D21').Value = '0.0₃0828'
$ws.Range('E21').Value = '  -1.12%  '
$ws.Range('D22').Value = '226.17'
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').Value = '2.29'
$ws.Range('E25').Value = '  -4.19%  '
$ws.Range('D26').Value = '9.70'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').Value = '168.09'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('D28').Value = '0.129'
$ws.Range('E28').Value = '  -5.99%  '
$ws.Range('D29').Value = '18.99'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('D30').Value = '1.36'
$ws.Range('E30').Value = '  -2.89%  '
$ws.Range('E31').Value = '  -1.66%  '
$ws.Range('D32').Value = '4.54'
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('D33').Value = '0.0615'
$ws.Range('E33').Value = '  -2.44%  '
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('E36').Value = '  +1.50%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('D38').Value = '3.25'
$ws.Range('E38').Value = '  -3.82%  '
$ws.Range('D39').Value = '5.26'
$ws.Range('E39').Value = '  -2.48%  '
$ws.Range('D40').Value = '0.0221'
$ws.Range('E40').Value = '  -4.29%  '
$ws.Range('D41').Value = '17.18'
$ws.Range('E41').Value = '  +1.23%  '
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').Value = '1.475.51'
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').Value = '0.0945'
$ws.Range('E44').Value = '  -3.02%  '
$ws.Range('D45').Value = '96.42'
$ws.Range('E45').Value = '  -4.68%  '
$ws.Range('E46').Value = '  +1.23%  '
$ws.Range('E47').Value = '  -3.71%  '
$ws.Range('D48').Value = '3.94'
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('D49').Value = '7.13'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('D51').Value = '2.228.83'
$ws.Range('E51').Value = '  -1.71%  '

Write-Host "Applied $(86) changes"
